$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column (D) keeps its text formatting so that
# numeric-looking strings (e.g. "1.00", "0.0864") are not converted
# into real numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '52.246.66'
$ws.Range("E2").Value = '  +0.59%  '
$ws.Range("D3").Value = '2.915.00'
$ws.Range("E3").Value = '  +3.48%  '
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").Value = '351.87'
$ws.Range("E5").Value = '  -1.69%  '
$ws.Range("D6").Value = '112.45'
$ws.Range("E6").Value = '  +2.19%  '
$ws.Range("D7").Value = '0.559'
$ws.Range("E7").Value = '  +0.11%  '
$ws.Range("E8").Value = '  +0.02%  '
$ws.Range("D9").Value = '0.630'
$ws.Range("E9").Value = '  -1.09%  '
$ws.Range("D10").Value = '39.99'
$ws.Range("E10").Value = '  -0.58%  '
$ws.Range("D11").Value = '0.0864'
$ws.Range("E11").Value = '  +2.56%  '
$ws.Range("D12").Value = '0.136'
$ws.Range("E12").Value = '  +0.16%  '
$ws.Range("D13").Value = '19.97'
$ws.Range("E13").Value = '  -0.55%  '
$ws.Range("D14").Value = '7.80'
$ws.Range("E14").Value = '  -0.35%  '
$ws.Range("D15").Value = '3.371.46'
$ws.Range("E15").Value = '  +3.50%  '
$ws.Range("B16").Value = 'WrappedEther'
$ws.Range("C16").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D16").Value = '2.910.96'
$ws.Range("E16").Value = '  +2.87%  '
$ws.Range("B17").Value = 'Polygon'
$ws.Range("C17").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D17").Value = '1.00'
$ws.Range("E17").Value = '  +5.89%  '
$ws.Range("D18").Value = '52.260.46'
$ws.Range("E18").Value = '  +0.71%  '
$ws.Range("D19").Value = '7.65'
$ws.Range("E19").Value = '  -0.74%  '
$ws.Range("D20").Value = '3.32'
$ws.Range("E20").Value = '  +4.19%  '
$ws.Range("D21").Value = '14.20'
$ws.Range("E21").Value = '  +3.72%  '
$ws.Range("D22").Value = '0.0₃0980'
$ws.Range("E22").Value = '  -0.07%  '
$ws.Range("D23").Value = '70.82'
$ws.Range("E23").Value = '  +0.42%  '
$ws.Range("D24").Value = '270.27'
$ws.Range("E24").Value = '  +0.44%  '
$ws.Range("E25").Value = '  +0.76%  '
$ws.Range("D26").Value = '26.78'
$ws.Range("E26").Value = '  +2.00%  '
$ws.Range("D28").Value = '0.166'
$ws.Range("E28").Value = '  +0.86%  '
$ws.Range("D29").Value = '10.62'
$ws.Range("E29").Value = '  +2.03%  '
$ws.Range("D30").Value = '37.43'
$ws.Range("E30").Value = '  -2.09%  '
$ws.Range("D31").Value = '2.25'
$ws.Range("E31").Value = '  +0.51%  '
$ws.Range("D32").Value = '6.47'
$ws.Range("E32").Value = '  +4.31%  '
$ws.Range("D33").Value = '6.13'
$ws.Range("E33").Value = '  +7.57%  '
$ws.Range("D34").Value = '0.0962'
$ws.Range("E34").Value = '  +10.56%  '
$ws.Range("D35").Value = '53.12'
$ws.Range("E35").Value = '  +1.83%  '
$ws.Range("D36").Value = '0.0452'
$ws.Range("E36").Value = '  +0.95%  '
$ws.Range("E37").Value = '  -0.07%  '
$ws.Range("D38").Value = '3.31'
$ws.Range("E38").Value = '  +4.93%  '
$ws.Range("D39").Value = '18.84'
$ws.Range("E39").Value = '  -0.47%  '
$ws.Range("D40").Value = '2.07'
$ws.Range("E40").Value = '  +2.77%  '
$ws.Range("D41").Value = '2.85'
$ws.Range("E41").Value = '  +12.98%  '
$ws.Range("B42").Value = 'EnergySwap'
$ws.Range("C42").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D42").Value = '23.63'
$ws.Range("E42").Value = '  +7.20%  '
$ws.Range("B43").Value = 'Stellar'
$ws.Range("C43").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D43").Value = '0.117'
$ws.Range("E43").Value = '  +1.08%  '
$ws.Range("E44").Value = '  +5.43%  '
$ws.Range("D45").Value = '120.66'
$ws.Range("E45").Value = '  +1.07%  '
$ws.Range("E46").Value = '  -0.33%  '
$ws.Range("B47").Value = 'NEARProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D47").Value = '3.55'
$ws.Range("E47").Value = '  +4.29%  '
$ws.Range("B48").Value = 'Maker'
$ws.Range("C48").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D48").Value = '2.196.01'
$ws.Range("E48").Value = '  +3.92%  '
$ws.Range("E49").Value = '  +23.01%  '
$ws.Range("D50").Value = '0.964'
$ws.Range("E50").Value = '  +3.36%  '
$ws.Range("D51").Value = '0.0335'
$ws.Range("E51").Value = '  +10.24%  '
